$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: row 31
$ws.Range("A31").Value = Get-Date -Year 2012 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("B31").Value = 2.5
$ws.Range("C31").Value = 1.25
$ws.Range("D31").Value = "Makefile: Dependency files incorporated, Manual continued"

# Match the date format/style used by the other cells in column A
$ws.Range("A31").NumberFormat = "ddd\ dd/mm/yyyy"

# Update selection to match the diff (activeCell A32, sqref A32)
$ws.Range("A32").Select()
